$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8").Value = "24h 30m"
$ws.Range("B8").Select()
